$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The league database refresh reordered several match records whose rows
# ended up adjacent to each other; for each such pair, row N's data
# (columns B:AC -- everything except the running index in column A) needs
# to be swapped with row N+1's data.
$pairs = @(
    @(742, 743),
    @(776, 777),
    @(778, 779),
    @(832, 833),
    @(866, 867),
    @(946, 947),
    @(957, 958)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:AC$r1")
    $range2 = $ws.Range("B$r2`:AC$r2")

    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value2 = $values2
    $range2.Value2 = $values1
}
